$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42:56 down to 43:57.
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new record.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 45029
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107001
$ws.Range("J42").Value = "Caqui"
$ws.Range("K42").Value = "Fuyu"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 130
$ws.Range("N42").Value = 22000
$ws.Range("O42").Value = 24000
$ws.Range("P42").Value = 22769
$ws.Range("Q42").Value = "$/bandeja 15 kilos granel"
$ws.Range("R42").Value = "Región de O'Higgins"
$ws.Range("S42").Value = 1518
$ws.Range("T42").Value = 15
